$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# This sheet holds weekly price records, most recent first (around row 82).
# A new weekly record is being added: insert a new row above row 82, which
# pushes the former row 82 down to row 83 and the former row 83 down to row 84.
$ws.Rows.Item(82).Insert()

# Fill the new row 82 with the new weekly record. Most fields repeat the
# values of the record directly below (now row 83), only the date and the
# price/volume figures change.
$ws.Range("A82").Value = 4
$ws.Range("B82").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C82").Value = "Los Lagos"
$ws.Range("D82").Value = 44595
$ws.Range("D82").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E82").Value = 10
$ws.Range("F82").Value = 100112026
$ws.Range("G82").Value = "Haba"
$ws.Range("H82").Value = "Sin especificar"
$ws.Range("I82").Value = "Primera"
$ws.Range("J82").Value = 40
$ws.Range("K82").Value = 26000
$ws.Range("L82").Value = 26000
$ws.Range("M82").Value = 26000
$ws.Range("N82").Value = "$/saco 25 kilos"
$ws.Range("O82").Value = "Región de La Araucanía"
$ws.Range("P82").Value = 1040
$ws.Range("Q82").Value = 25
$ws.Range("R82").Value = "Hortaliza"
